$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '31.090.10'
$ws.Range('E2').Value = '  +4.12%  '
$ws.Range('D3').Value = '1.683.48'
$ws.Range('E3').Value = '  +3.26%  '
$ws.Range('E4').Value = '  -0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '220.35'
$ws.Range('E5').Value = '  +2.60%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.534'
$ws.Range('E6').Value = '  +2.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '29.48'
$ws.Range('E8').Value = '  +2.74%  '
$ws.Range('E9').Value = '  +3.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0640'
$ws.Range('E10').Value = '  +5.18%  '
$ws.Range('E11').Value = '  +0.97%  '
$ws.Range('D12').Value = '1.922.80'
$ws.Range('E12').Value = '  +3.21%  '
$ws.Range('D13').Value = '1.682.57'
$ws.Range('E13').Value = '  +3.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.35'
$ws.Range('E14').Value = '  +9.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.609'
$ws.Range('E15').Value = '  +7.18%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.10'
$ws.Range('E16').Value = '  +7.11%  '
$ws.Range('D17').Value = '31.059.46'
$ws.Range('E17').Value = '  +3.97%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '66.55'
$ws.Range('E18').Value = '  +2.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '248.37'
$ws.Range('E19').Value = '  +3.25%  '
$ws.Range('D20').Value = '0.0₃0721'
$ws.Range('E20').Value = '  +2.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.28'
$ws.Range('E22').Value = '  +3.60%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.01'
$ws.Range('E23').Value = '  +2.25%  '
$ws.Range('E24').Value = '  -0.83%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '158.79'
$ws.Range('E25').Value = '  +0.79%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.92'
$ws.Range('E26').Value = '  +2.90%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.113'
$ws.Range('E27').Value = '  +2.79%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.71'
$ws.Range('E28').Value = '  +1.63%  '
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('E30').Value = '  +1.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.52'
$ws.Range('E31').Value = '  +4.15%  '
$ws.Range('E32').Value = '  +3.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.34'
$ws.Range('E33').Value = '  +5.24%  '
$ws.Range('D34').Value = '1.516.88'
$ws.Range('E34').Value = '  +6.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.74'
$ws.Range('E35').Value = '  +3.70%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '84.80'
$ws.Range('E36').Value = '  +13.27%  '
$ws.Range('E37').Value = '  +0.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.613'
$ws.Range('E38').Value = '  +10.37%  '
$ws.Range('E39').Value = '  +5.44%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.67'
$ws.Range('E40').Value = '  -3.46%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.05'
$ws.Range('E42').Value = '  +3.40%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.843'
$ws.Range('E43').Value = '  +1.44%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0504'
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('E45').Value = '  +2.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.998'
$ws.Range('E46').Value = '  -0.24%  '
$ws.Range('B47').Value = 'BitcoinSV'
$ws.Range('C47').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '52.10'
$ws.Range('E47').Value = '  +6.66%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.58'
$ws.Range('E48').Value = '  +4.72%  '
$ws.Range('E49').Value = '  +2.33%  '
$ws.Range('D50').Value = '0.0₆0120'
$ws.Range('E50').Value = '  +8.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '93.51'
$ws.Range('E51').Value = '  +1.45%  '
